$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Move the <w:lastRenderedPageBreak/> marker from the "Start
#    results" paragraph to the "Step one: assess progress..." one.
# -----------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Start results", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $rng.Paragraphs(1)
$para.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Start results</w:t></w:r></w:p>')

$rng = $d.Content
$rng.Find.Execute("Step one: assess progress so far and skim through old code", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $rng.Paragraphs(1)
$para.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Step one: assess progress so far and skim through old code</w:t></w:r></w:p>')

# -----------------------------------------------------------------
# 2) Rework the BcBOT bigRR / cis-trans effects outline: the block
#    that used to run from "BcBOT bigRR" through "Not sure where
#    list came from/ how to verify" is replaced wholesale with the
#    restructured, expanded outline (new "Troubleshoot small cis
#    effects" heading, the crude-haplotype / fastPhase / trans
#    effects sub-outline, and the old sub-items reattached below).
# -----------------------------------------------------------------
$startRng = $d.Content
$startRng.Find.Execute("BcBOT bigRR", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPara = $startRng.Paragraphs(1)

$endRng = $d.Content
$endRng.Find.Execute("Are low expression isolates genetically related in this gene?", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPara = $endRng.Paragraphs(1)

$blockRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$blockRange.InsertXML('    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>Troubleshoot small cis effects</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:rPr>
          <w:color w:val="808080" w:themeColor="background1" w:themeShade="80"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Select 1 </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>BcBOT</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> gene</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:rPr>
          <w:color w:val="808080" w:themeColor="background1" w:themeShade="80"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Dendrogram</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>: use all SNPs within gene ends (+- 1 kb?)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:rPr>
          <w:color w:val="808080" w:themeColor="background1" w:themeShade="80"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t>Are low expression isolates genetically related in this gene?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve"> Crude haplotype calls per gene based on SNPs within gene</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Problem: SNP by SNP association does not capture major cis-effects (due to minor allele/ haplotype)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Approach: lasso on SNP correlation matrix for locus?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Step one: crude haplotype assay to grab cis effects</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>fastPhase</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">: see scripts from Michelle </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Stitzer</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="4"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>sensitive to a priori number of haplotypes</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="4"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>grab haplotype assignment surrounding a focal site or break into longest segments</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Step two: map residuals after cis effects </w:t>
      </w:r>
      <w:r>
        <w:sym w:font="Wingdings" w:char="F0E0"/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> ID trans effects</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>BcBOT</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>bigRR</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Quality check work so far</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Check for completeness of </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>BcBOT</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> list previous run</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:rPr>
          <w:color w:val="808080" w:themeColor="background1" w:themeShade="80"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="808080" w:themeColor="background1" w:themeShade="80"/>
        </w:rPr>
        <w:t>Not sure where list came from/ how to verify</w:t>
      </w:r>
    </w:p>')
